# Insert a new "pander(table_forecasts)" SourceCode paragraph right after
# the "Using historical data..." BodyText paragraph and before the table
# of McDonald's stock forecasts that follows it.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "...summarized in the table below:"
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -like "*Using historical data*summarized in the table below:*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq 0) {
    throw "Could not locate the 'Using historical data' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Create a new, empty paragraph right after it (before the table).
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Give the new paragraph the "SourceCode" paragraph style.
$newPara.Style = "SourceCode"

# Insert the literal code text as plain text first.
$newPara.Range.InsertAfter("pander(table_forecasts)")

# Now split it into two runs by re-styling each piece via Find/Replace
# (this applies a character style without disturbing the rest of the
# document).
$r1 = $newPara.Range
$f1 = $r1.Find
$f1.ClearFormatting()
$f1.Replacement.ClearFormatting()
$f1.Replacement.Style = "FunctionTok"
[void]$r1.Find.Execute("pander", $true, $false, $false, $false, $false, $true, 1, $true, "^&", 2)

$r2 = $newPara.Range
$f2 = $r2.Find
$f2.ClearFormatting()
$f2.Replacement.ClearFormatting()
$f2.Replacement.Style = "NormalTok"
[void]$r2.Find.Execute("(table_forecasts)", $true, $false, $false, $false, $false, $true, 1, $true, "^&", 2)

Write-Output "New paragraph text: $($newPara.Range.Text)"
